# Update crypto price/volume figures (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.533.07"
$ws.Range("E2").Value = "  +2.31%  "

$ws.Range("D3").Value = "3.459.59"
$ws.Range("E3").Value = "  +2.63%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "'575.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.48%  "

$ws.Range("D6").Value = "'158.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.58%  "

$ws.Range("E7").Value = "  -0.15%  "

$ws.Range("D8").Value = "3.468.17"
$ws.Range("E8").Value = "  +2.66%  "

$ws.Range("D9").Value = "'0.581"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +10.34%  "

$ws.Range("D10").Value = "'7.37"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.86%  "

$ws.Range("E11").Value = "  +5.30%  "

$ws.Range("D12").Value = "'0.444"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.98%  "

$ws.Range("D13").Value = "4.043.91"
$ws.Range("E13").Value = "  +2.33%  "

$ws.Range("E14").Value = "  -2.82%  "

$ws.Range("E15").Value = "  +7.38%  "

$ws.Range("D16").Value = "'28.31"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.02%  "

$ws.Range("D17").Value = "64.480.74"
$ws.Range("E17").Value = "  +2.12%  "

$ws.Range("D18").Value = "3.447.29"
$ws.Range("E18").Value = "  +3.46%  "

$ws.Range("D19").Value = "'6.45"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.13%  "

$ws.Range("D20").Value = "'14.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.83%  "

$ws.Range("D21").Value = "'390.50"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.33%  "

$ws.Range("D22").Value = "'8.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.60%  "

$ws.Range("D23").Value = "'73.73"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.71%  "

$ws.Range("D24").Value = "'0.544"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.03%  "

$ws.Range("E25").Value = "  +0.01%  "

$ws.Range("E26").Value = "  +24.66%  "

$ws.Range("D27").Value = "'9.60"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.47%  "

$ws.Range("E28").Value = "  +0.49%  "

$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.10%  "

$ws.Range("D30").Value = "'6.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +11.14%  "

$ws.Range("E31").Value = "  +10.49%  "

$ws.Range("E32").Value = "  +0.74%  "

$ws.Range("D33").Value = "'6.55"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.88%  "

$ws.Range("E34").Value = "  +2.77%  "

$ws.Range("D35").Value = "'0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.07%  "

$ws.Range("D36").Value = "'7.03"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.83%  "

$ws.Range("D37").Value = "'1.49"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.81%  "

$ws.Range("D38").Value = "'160.81"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.58%  "

$ws.Range("E39").Value = "  +1.69%  "

$ws.Range("D40").Value = "'0.0775"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.46%  "

$ws.Range("D41").Value = "2.930.66"
$ws.Range("E41").Value = "  +1.26%  "

$ws.Range("D42").Value = "'27.33"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.06%  "

$ws.Range("D43").Value = "'0.0319"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.45%  "

$ws.Range("D44").Value = "'42.69"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.61%  "

$ws.Range("D45").Value = "'4.44"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.96%  "

$ws.Range("D46").Value = "'0.772"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.08%  "

$ws.Range("D47").Value = "'23.76"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +8.44%  "

$ws.Range("D48").Value = "'1.09"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.89%  "

$ws.Range("D49").Value = "'2.25"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +22.68%  "

$ws.Range("E50").Value = "  +4.66%  "

$ws.Range("D51").Value = "'0.864"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.43%  "
